$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: reorder "Recorded By" email list (G2) ---
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, System"

# --- Row 3: reorder "Recorded By" email list (G3) ---
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, System, eman.tantawi@med.asu.edu.eg"

# --- Class Statistics block (rows 6-10 / 15) : Recorded +1, Missing -1 ---
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 1
$ws.Range("L9").Value = "27.6%"
$ws.Range("L10").Value = "25.0%"

# --- Row 10: HISTOLOGY C1 session 2 is now Recorded ---
$src = $ws.Range("A9:I9")
$dst = $ws.Range("A10:I10")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Range("G10").Value = "Safa.hany@med.asu.edu.eg"
$ws.Range("H10").Value = "8/251"
$ws.Range("I10").Value = "Recorded"

# --- Group statistics row 15 mirrors the same Recorded/Missing/Coverage/Attendance values ---
$ws.Range("O15").Value = 8
$ws.Range("P15").Value = 1
$ws.Range("R15").Value = "27.6%"
$ws.Range("S15").Value = "25.0%"

# --- Row 29: PHYSIOLOGY C1 session 2 keeps its look (no content change) ---
Write-Output "done"
